$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.787066666666667
$ws.Range("H2").Value = 5.3612
$ws.Range("I2").Value = 0.4880702634734595
$ws.Range("J2").Value = 0.4880702634734594
$ws.Range("O2").Value = 0.00762120683641437
$ws.Range("P2").Value = 0.007621206836414371
$ws.Range("Q2").Value = 0.0002603160444444445
$ws.Range("R2").Value = 0.0023428444
$ws.Range("S2").Value = 0.003719684428634492
$ws.Range("T2").Value = 0.003719684428634492

$ws.Range("G3").Value = 1.787066666666667
$ws.Range("H3").Value = 5.3612
$ws.Range("I3").Value = 0.4880702634734595
$ws.Range("J3").Value = 0.4880702634734594
$ws.Range("O3").Value = 0.05683641437042205
$ws.Range("P3").Value = 0.05683641437042205
$ws.Range("Q3").Value = 0.001941350088888889
$ws.Range("R3").Value = 0.0174721508
$ws.Range("S3").Value = 0.02774016373665861
$ws.Range("T3").Value = 0.0277401637366586

$ws.Range("G4").Value = 1.787066666666667
$ws.Range("H4").Value = 5.3612
$ws.Range("I4").Value = 0.4880702634734595
$ws.Range("J4").Value = 0.4880702634734594
$ws.Range("M4").Value = 0.01788133333333333
$ws.Range("Q4").Value = 0.03195513475555556
$ws.Range("R4").Value = 0.2875962128
$ws.Range("S4").Value = 0.4566104153081664
$ws.Range("T4").Value = 0.4566104153081663

$ws.Range("G5").Value = 1.391264666666667
$ws.Range("H5").Value = 4.173794
$ws.Range("I5").Value = 0.3799717856569321
$ws.Range("J5").Value = 0.379971785656932
$ws.Range("O5").Value = 0.00762120683641437
$ws.Range("P5").Value = 0.007621206836414371
$ws.Range("S5").Value = 0.002895843570493186
$ws.Range("T5").Value = 0.002895843570493186

$ws.Range("G6").Value = 1.391264666666667
$ws.Range("H6").Value = 4.173794
$ws.Range("I6").Value = 0.3799717856569321
$ws.Range("J6").Value = 0.379971785656932
$ws.Range("O6").Value = 0.05683641437042205
$ws.Range("P6").Value = 0.05683641437042205
$ws.Range("S6").Value = 0.02159623385866658
$ws.Range("T6").Value = 0.02159623385866658

$ws.Range("G7").Value = 1.391264666666667
$ws.Range("H7").Value = 4.173794
$ws.Range("I7").Value = 0.3799717856569321
$ws.Range("J7").Value = 0.379971785656932
$ws.Range("M7").Value = 0.01788133333333333
$ws.Range("S7").Value = 0.3554797082277723
$ws.Range("T7").Value = 0.3554797082277723

$ws.Range("G8").Value = 0.159805
$ws.Range("H8").Value = 0.479415
$ws.Range("I8").Value = 0.04364474471445359
$ws.Range("J8").Value = 0.04364474471445358
$ws.Range("O8").Value = 0.00762120683641437
$ws.Range("P8").Value = 0.007621206836414371
$ws.Range("Q8").Value = 0.00002327826166666667
$ws.Range("R8").Value = 0.000209504355
$ws.Range("S8").Value = 0.0003326256267913536
$ws.Range("T8").Value = 0.0003326256267913536

$ws.Range("G9").Value = 0.159805
$ws.Range("H9").Value = 0.479415
$ws.Range("I9").Value = 0.04364474471445359
$ws.Range("J9").Value = 0.04364474471445358
$ws.Range("O9").Value = 0.05683641437042205
$ws.Range("P9").Value = 0.05683641437042205
$ws.Range("Q9").Value = 0.0001736014983333333
$ws.Range("R9").Value = 0.001562413485
$ws.Range("S9").Value = 0.002480610795681972
$ws.Range("T9").Value = 0.002480610795681971

$ws.Range("G10").Value = 0.159805
$ws.Range("H10").Value = 0.479415
$ws.Range("I10").Value = 0.04364474471445359
$ws.Range("J10").Value = 0.04364474471445358
$ws.Range("M10").Value = 0.01788133333333333
$ws.Range("Q10").Value = 0.002857526473333333
$ws.Range("R10").Value = 0.02571773826
$ws.Range("S10").Value = 0.04083150829198027
$ws.Range("T10").Value = 0.04083150829198025

$ws.Range("G11").Value = 0.1755086666666667
$ws.Range("H11").Value = 0.526526
$ws.Range("I11").Value = 0.04793361253928724
$ws.Range("J11").Value = 0.04793361253928723
$ws.Range("O11").Value = 0.00762120683641437
$ws.Range("P11").Value = 0.007621206836414371
$ws.Range("Q11").Value = 0.00002556576244444445
$ws.Range("R11").Value = 0.000230091862
$ws.Range("S11").Value = 0.0003653119755784535
$ws.Range("T11").Value = 0.0003653119755784535

$ws.Range("G12").Value = 0.1755086666666667
$ws.Range("H12").Value = 0.526526
$ws.Range("I12").Value = 0.04793361253928724
$ws.Range("J12").Value = 0.04793361253928723
$ws.Range("O12").Value = 0.05683641437042205
$ws.Range("P12").Value = 0.05683641437042205
$ws.Range("Q12").Value = 0.0001906609148888889
$ws.Range("R12").Value = 0.001715948234
$ws.Range("S12").Value = 0.002724374664554188
$ws.Range("T12").Value = 0.002724374664554188

$ws.Range("G13").Value = 0.1755086666666667
$ws.Range("H13").Value = 0.526526
$ws.Range("I13").Value = 0.04793361253928724
$ws.Range("J13").Value = 0.04793361253928723
$ws.Range("M13").Value = 0.01788133333333333
$ws.Range("Q13").Value = 0.003138328971555555
$ws.Range("R13").Value = 0.028244960744
$ws.Range("S13").Value = 0.0448439258991546
$ws.Range("T13").Value = 0.04484392589915459

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1478496666666667
$ws.Range("H14").Value = 0.443549
$ws.Range("I14").Value = 0.04037959361586762
$ws.Range("J14").Value = 0.04037959361586761
$ws.Range("O14").Value = 0.00762120683641437
$ws.Range("P14").Value = 0.007621206836414371
$ws.Range("Q14").Value = 0.00002153676811111111
$ws.Range("R14").Value = 0.000193830913
$ws.Range("S14").Value = 0.0003077412349168843
$ws.Range("T14").Value = 0.0003077412349168843

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1478496666666667
$ws.Range("H15").Value = 0.443549
$ws.Range("I15").Value = 0.04037959361586762
$ws.Range("J15").Value = 0.04037959361586761
$ws.Range("O15").Value = 0.05683641437042205
$ws.Range("P15").Value = 0.05683641437042205
$ws.Range("Q15").Value = 0.0001606140212222222
$ws.Range("R15").Value = 0.001445526191
$ws.Range("S15").Value = 0.002295031314860701
$ws.Range("T15").Value = 0.0022950313148607

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1478496666666667
$ws.Range("H16").Value = 0.443549
$ws.Range("I16").Value = 0.04037959361586762
$ws.Range("J16").Value = 0.04037959361586761
$ws.Range("M16").Value = 0.01788133333333333
$ws.Range("Q16").Value = 0.002643749172888888
$ws.Range("R16").Value = 0.023793742556
$ws.Range("S16").Value = 0.03777682106609003
$ws.Range("T16").Value = 0.03777682106609002
